$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 309
$ws1.Range("F4").Value = 1249

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 309
$ws4.Range("F4").Value = 1249
